$wb = $excel.ActiveWorkbook
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D6").Value = "2016-01-17 03:14:38"
$wsDe.Range("D6").Value = "2016-01-17 03:14:47"
